{"js": "// The diff removes three consecutive paragraphs that used to sit right\n// after the \"LOM3058: ...\" requirement line at the end of the document:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//       pages. Original theme under Creative Commons Attribution\"\n// The empty paragraph and page-break paragraph that originally followed\n// those three paragraphs are left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraphs to delete by their text content so the script is\n// resilient to the exact paragraph indices.\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nlet jupiterIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === jupiterText) {\n    jupiterIdx = i;\n    break;\n  }\n}\n\nif (jupiterIdx === -1) {\n  throw new Error('Could not find paragraph \"' + jupiterText + '\"');\n}\n\nconst emptyIdx = jupiterIdx - 1;\nconst copyrightIdx = jupiterIdx + 1;\n\nif (items[copyrightIdx].text !== copyrightText) {\n  throw new Error(\"Unexpected paragraph following the Jupiter paragraph\");\n}\nif (items[emptyIdx].text !== \"\") {\n  throw new Error(\"Unexpected paragraph preceding the Jupiter paragraph\");\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIdx].delete();\nitems[jupiterIdx].delete();\nitems[emptyIdx].delete();\n\nawait context.sync();\n", "ps1": "# The diff removes three consecutive paragraphs that used to sit right\n# after the \"LOM3058: ...\" requirement line at the end of the document:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) \"(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#       pages. Original theme under Creative Commons Attribution\"\n# The empty paragraph and page-break paragraph that originally followed\n# those three paragraphs are left untouched.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"LOM3058*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the LOM3058 paragraph\"\n}\n\n$p1 = $target.Next()   # empty paragraph\n$p2 = $p1.Next()       # \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$p3 = $p2.Next()       # \"(c) 2020 . Contact: ... Creative Commons Attribution\"\n\nif ($p2.Range.Text -notlike \"Ver no Jupiter*\") {\n    throw \"Unexpected paragraph where 'Ver no Jupiter...' was expected\"\n}\nif ($p3.Range.Text -notlike \"*Creative Commons*\") {\n    throw \"Unexpected paragraph where the copyright line was expected\"\n}\n\n# Build a single range spanning from the start of the empty paragraph to\n# the end of the copyright paragraph (Range.End already covers its\n# paragraph mark) and delete it in one shot.\n$delRange = $d.Range($p1.Range.Start, $p3.Range.End)\n$delRange.Delete()\n"}
